$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing AgTests (F) / AgPosit (G) values ---
$ws.Range("F668").Value2 = 3411

$ws.Range("F715").Value2 = 32158

$ws.Range("F740").Value2 = 25348

$ws.Range("F756").Value2 = 13878
$ws.Range("G756").Value2 = 1077

$ws.Range("F757").Value2 = 13672

$ws.Range("F779").Value2 = 7388

$ws.Range("F782").Value2 = 10729

$ws.Range("F783").Value2 = 7845

$ws.Range("F784").Value2 = 7816

$ws.Range("F785").Value2 = 7185

$ws.Range("F786").Value2 = 6341

$ws.Range("F789").Value2 = 7880

$ws.Range("F791").Value2 = 4252

$ws.Range("F793").Value2 = 3407

$ws.Range("F796").Value2 = 4656

$ws.Range("F798").Value2 = 3509

$ws.Range("F800").Value2 = 3078

$ws.Range("F803").Value2 = 3821

$ws.Range("F805").Value2 = 2408

$ws.Range("F807").Value2 = 2327

$ws.Range("F810").Value2 = 3894

$ws.Range("F812").Value2 = 2102

$ws.Range("F814").Value2 = 1959

$ws.Range("F817").Value2 = 3762

$ws.Range("F818").Value2 = 2198

$ws.Range("F819").Value2 = 2259

$ws.Range("F821").Value2 = 1654

$ws.Range("F824").Value2 = 3589
$ws.Range("G824").Value2 = 95

$ws.Range("F825").Value2 = 4430
$ws.Range("G825").Value2 = 48

$ws.Range("F826").Value2 = 1824
$ws.Range("G826").Value2 = 65

# --- Row 827 was missing F/G values; fill them in ---
$ws.Range("F827").Value2 = 1899
$ws.Range("G827").Value2 = 39

# --- Append new daily rows 828-833 ---
$ws.Range("A828").Value2 = 44722
$ws.Range("B828").Value2 = 1791124
$ws.Range("C828").Value2 = 1084
$ws.Range("D828").Value2 = 165
$ws.Range("E828").Value2 = 20114
$ws.Range("F828").Value2 = 1289
$ws.Range("G828").Value2 = 46

$ws.Range("A829").Value2 = 44723
$ws.Range("B829").Value2 = 1791192
$ws.Range("C829").Value2 = 524
$ws.Range("D829").Value2 = 68
$ws.Range("E829").Value2 = 20116
$ws.Range("F829").Value2 = 531
$ws.Range("G829").Value2 = 18

$ws.Range("A830").Value2 = 44724
$ws.Range("B830").Value2 = 1791220
$ws.Range("C830").Value2 = 244
$ws.Range("D830").Value2 = 28
$ws.Range("E830").Value2 = 20116
$ws.Range("F830").Value2 = 620
$ws.Range("G830").Value2 = 37

$ws.Range("A831").Value2 = 44725
$ws.Range("B831").Value2 = 1791426
$ws.Range("C831").Value2 = 1265
$ws.Range("D831").Value2 = 206
$ws.Range("E831").Value2 = 20119
$ws.Range("F831").Value2 = 1961
$ws.Range("G831").Value2 = 104

$ws.Range("A832").Value2 = 44726
$ws.Range("B832").Value2 = 1791639
$ws.Range("C832").Value2 = 1032
$ws.Range("D832").Value2 = 213
$ws.Range("E832").Value2 = 20119
$ws.Range("F832").Value2 = 1823
$ws.Range("G832").Value2 = 54

$ws.Range("A833").Value2 = 44727
$ws.Range("B833").Value2 = 1791816
$ws.Range("C833").Value2 = 918
$ws.Range("D833").Value2 = 177
$ws.Range("E833").Value2 = 20122
$ws.Range("F833").Value2 = 747
$ws.Range("G833").Value2 = 44

# Ensure the new date cells (column A) use the same date style as the rest of the column
$ws.Range("A828:A833").NumberFormat = "yyyy-mm-dd"
